$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MainSheet")

# --- Cell value updates (rows shifted: Run/No row moved, device list refreshed) ---
$ws.Range("A2").Value = "Run"
$ws.Range("B2").Value = "webApp"
$ws.Range("C2").Value = "Web"
$ws.Range("D2").Value = "Windows_Chrome"
$ws.Range("E2").Value = "None"

$ws.Range("A3").Value = "No"
$ws.Range("D3").Value = "Android_6.1"
$ws.Range("E3").Value = "ZY32288VFB"
$ws.Range("F3").Value = "No"

$ws.Range("F4").Value = "No"
$ws.Range("H4").Value = "Yes"

# --- Selection moves to A4 ---
$ws.Range("A4").Select()

# --- Data validation updates ---

# Remove the old E2:E4 validation now so the list re-order below ends with it last.
$ws.Range("E2:E4").Validation.Delete()

# Fix the AppType (B2) validation list casing: NativeApp -> nativeApp.
$ws.Range("B2").Validation.Delete()
$ws.Range("B2").Validation.Add(3, 1, 1, """webApp,nativeApp,sanityTesting""")
$ws.Range("B2").Validation.IgnoreBlank = $false

# Add Android_5.0 option to the Device Name (D2:D4) validation list.
$ws.Range("D2:D4").Validation.Delete()
$ws.Range("D2:D4").Validation.Add(3, 1, 1, """Windows_Chrome,Windows_Firefox,Windows_Safari,Android_6.0,Android_7.1,Android_5.0""")

# Re-add the UDID (E2:E4) validation with the new device ids, now at the end of the list.
$ws.Range("E2:E4").Validation.Add(3, 1, 1, """None,emulator-5554,emulator-5556,b1a1589f,ZY32288VFB""")
$ws.Range("E2:E4").Validation.IgnoreBlank = $false
